$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

$text = "數位相框`nDigital Frame`nLoading Bitmaps...`nPlease wait for a moment.`n載入中..."

$ws.Range("B33").Value = "SingleUseId35"
$ws.Range("C33").Value = "Chinese"
$ws.Range("D33").Value = "Center"
$ws.Range("E33").Value = "LTR"
$ws.Range("F33").Value = $text

$ws.Rows("33").AutoFit()
